$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "平潭发展"
$ws.Range("A3").Value = "合富中国"
$ws.Range("A4").Value = "孚日股份"
$ws.Range("B4").Value = "海马汽车"
$ws.Range("C4").Value = "海马汽车"
$ws.Range("A5").Value = "华夏幸福"
$ws.Range("B5").Value = "多氟多"
$ws.Range("C5").Value = "安泰集团"
$ws.Range("A6").Value = "海马汽车"
$ws.Range("B6").Value = "华夏幸福"
$ws.Range("C6").Value = "孚日股份"
$ws.Range("A7").Value = "多氟多"
$ws.Range("B7").Value = "安泰集团"
$ws.Range("A8").Value = "天际股份"
$ws.Range("B8").Value = "孚日股份"
$ws.Range("C8").Value = "东百集团"
$ws.Range("B9").Value = "东百集团"
$ws.Range("C9").Value = "华夏幸福"
$ws.Range("A10").Value = "众生药业"
$ws.Range("B10").Value = "天赐材料"
$ws.Range("C10").Value = "多氟多"
$ws.Range("B11").Value = "隆基绿能"
$ws.Range("C11").Value = "粤桂股份"
$ws.Range("A12").Value = "人民同泰"
$ws.Range("B12").Value = "众生药业"
$ws.Range("C12").Value = "天际股份"
$ws.Range("A13").Value = "摩恩电气"
$ws.Range("B13").Value = "永泰能源"
$ws.Range("C13").Value = "三木集团"
$ws.Range("A14").Value = "隆基绿能"
$ws.Range("B14").Value = "摩恩电气"
$ws.Range("C14").Value = "永泰能源"
$ws.Range("A15").Value = "特一药业"
$ws.Range("B15").Value = "盈新发展"
$ws.Range("C15").Value = "国晟科技"
$ws.Range("A16").Value = "永泰能源"
$ws.Range("B16").Value = "康芝药业"
$ws.Range("C16").Value = "隆基绿能"
$ws.Range("A17").Value = "盈新发展"
$ws.Range("B17").Value = "海南海药"
$ws.Range("C17").Value = "工业富联"
$ws.Range("A18").Value = "福龙马"
$ws.Range("B18").Value = "特一药业"
$ws.Range("C18").Value = "特一药业"
$ws.Range("A19").Value = "天赐材料"
$ws.Range("B19").Value = "天际股份"
$ws.Range("C19").Value = "人民同泰"
$ws.Range("A20").Value = "三木集团"
$ws.Range("B20").Value = "海峡创新"
$ws.Range("C20").Value = "福龙马"
$ws.Range("A21").Value = "海峡创新"
$ws.Range("B21").Value = "福龙马"
$ws.Range("C21").Value = "盈新发展"
